$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C77").Value = "'55"
$ws.Range("D77").Value = "'201579.25"
$ws.Range("C80").Value = "'301"
$ws.Range("D80").Value = "'893226.16"
$ws.Range("C81").Value = "'67"
$ws.Range("D81").Value = "'193999.78"
$ws.Range("C82").Value = "'770"
$ws.Range("D82").Value = "'5691633.84"
$ws.Range("C91").Value = "'138"
$ws.Range("D91").Value = "'360280.06"
$ws.Range("C95").Value = "'259"
$ws.Range("D95").Value = "'825298.00"
$ws.Range("C96").Value = "'31"
$ws.Range("D96").Value = "'160500.00"
$ws.Range("C97").Value = "'247"
$ws.Range("D97").Value = "'1035230.00"
$ws.Range("C101").Value = "'77"
$ws.Range("D101").Value = "'368142.77"
$ws.Range("C102").Value = "'108"
$ws.Range("D102").Value = "'344423.04"
$ws.Range("C105").Value = "'42"
$ws.Range("D105").Value = "'247694.94"
$ws.Range("C136").Value = "'150"
$ws.Range("D136").Value = "'1020996.69"
$ws.Range("C144").Value = "'8514"
$ws.Range("D144").Value = "'27392204.19"
$ws.Range("C145").Value = "'5144"
$ws.Range("D145").Value = "'33582426.92"
$ws.Range("C148").Value = "'349"
$ws.Range("D148").Value = "'1199848.33"
$ws.Range("C150").Value = "'858"
$ws.Range("D150").Value = "'3254331.22"
$ws.Range("C153").Value = "'854"
$ws.Range("D153").Value = "'4311753.90"
$ws.Range("C256").Value = "'231"
$ws.Range("D256").Value = "'601290.00"
$ws.Range("C257").Value = "'768"
$ws.Range("D257").Value = "'2083212.76"
$ws.Range("C258").Value = "'527"
$ws.Range("D258").Value = "'1311632.81"
$ws.Range("C259").Value = "'1529"
$ws.Range("D259").Value = "'9131831.11"
$ws.Range("C262").Value = "'104"
$ws.Range("D262").Value = "'312442.28"
$ws.Range("C263").Value = "'267"
$ws.Range("D263").Value = "'1201366.83"
$ws.Range("C264").Value = "'194"
$ws.Range("D264").Value = "'829461.71"
$ws.Range("C265").Value = "'124"
$ws.Range("D265").Value = "'363705.92"
$ws.Range("C267").Value = "'246"
$ws.Range("D267").Value = "'1577391.53"
